$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells stay text so values like "62.796.97" are not
# reinterpreted as numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.796.97"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.464.46"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.60%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.36"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.17"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +17.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.457.52"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.698"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.17%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +29.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.06"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.71%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.015.11"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.81"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.12"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.458.36"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.832.62"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.45%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.81"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000140"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +27.29%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.70"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +10.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.19"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "313.24"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.17"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.45"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.99%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.01%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.68%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Hedera"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "InjectiveProtocol"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.26"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +12.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.79"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0493"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.65"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.59%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.99"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.54"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.52"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.289"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.97"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.15"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.810.11"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.180.74"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.12%  "
